$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix capitalization of "de"/"de la" -> "De"/"De La" in state/municipality names
$ws.Range("B4").Value  = "Amatenango De La Frontera"
$ws.Range("B7").Value  = "Chiapa De Corzo"
$ws.Range("A17").Value = "Ciudad De México"
$ws.Range("B20").Value = "Cuajimalpa De Morelos"
$ws.Range("A31").Value = "Estado De México"
$ws.Range("B38").Value = "Tlalnepantla De Baz"
$ws.Range("B43").Value = "Acapulco De Juárez"
$ws.Range("B46").Value = "Chilapa De Álvarez"
$ws.Range("B49").Value = "Iguala De La Independencia"
$ws.Range("B51").Value = "Tlapa De Comonfort"
$ws.Range("B53").Value = "Tulancingo De Bravo"
$ws.Range("B55").Value = "Lagos De Moreno"
$ws.Range("B62").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B70").Value = "Chalcatongo De Hidalgo"
$ws.Range("B71").Value = "Oaxaca De Juárez"
$ws.Range("B81").Value = "Izúcar De Matamoros"
$ws.Range("B87").Value = "Tepanco De López"
$ws.Range("B89").Value = "Tepexi De Rodríguez"
$ws.Range("B90").Value = "Tlacotepec De Benito Juárez"

# Remove trailing footer/notes rows (119-123), which sit below the data table (row 118 is blank/unused)
$ws.Rows("119:123").Delete()
